$wb = $excel.ActiveWorkbook

# --- Rename "kevin" sheet to "aboona" ---
$wsAboona = $wb.Worksheets.Item("kevin")
$wsAboona.Name = "aboona"

# --- "coach" sheet: bring the F (columns F:I) and K (columns K:N) team-group
#     header cells in rows 12-15 into line with the A/P groups' look
#     (drop the yellow highlight + mismatched borders, drop the manual
#     row-height/thick-bottom-border on the two interior rows) ---
$wsCoach = $wb.Worksheets.Item("coach")

$wsCoach.Range("A12:A15").Copy()
$wsCoach.Range("F12:F15").PasteSpecial(-4122)
$wsCoach.Range("K12:K15").PasteSpecial(-4122)

$wsCoach.Rows.Item(13).AutoFit()
$wsCoach.Rows.Item(14).AutoFit()

# --- restore the selections recorded for each sheet ---
$wsCoach.Activate()
$wsCoach.Range("I28").Select()

$wsAboona.Activate()
$wsAboona.Range("I35").Select()

Write-Output "done"
